# Add landing page translations (withdraw consent) to row 47.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings / row 47 values:
#   C47 -> en-US: "withdraw consent"
#   D47 -> de-DE: "Einwilligung wiederrufen"
#   E47 -> it-IT: "ritirare consenso"
#   F47 -> fr-FR: "retirer consentement"
$ws.Range("C47").Value = "withdraw consent"
$ws.Range("D47").Value = "Einwilligung wiederrufen"
$ws.Range("F47").Value = "retirer consentement"
$ws.Range("E47").Value = "ritirare consenso"

# Scroll/select so the active cell matches the new selection location.
$ws.Range("E47").Select()
